$wb = $excel.ActiveWorkbook

$wsWaterChem = $wb.Worksheets.Item("water chem")
$wsBenthic   = $wb.Worksheets.Item("raw benthic data")

# --- Data correction: rows 42-62 on "raw benthic data" used the wrong
# Activity_ID ("PRWI-MAWI"). Re-point them at "PRWI-MARU" so the site code
# matches the "water chem" / activities module, and pick up that sheet's
# cell formatting (style) for the corrected value while we're at it.
$wsWaterChem.Range("A4").Copy() | Out-Null
$targetRange = $wsBenthic.Range("A42:A62")
$targetRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$targetRange.Value = "PRWI-MARU"
$wsWaterChem.Application.CutCopyMode = $false

# --- View-state bookkeeping: the workbook was left with "raw benthic data"
# scrolled down near row 40 with B61 selected as the active tab, while
# "water chem" keeps a plain (non-active) selection on A4.
$wsWaterChem.Range("A4").Select() | Out-Null

$wsBenthic.Activate() | Out-Null
$benthicWindow = $excel.ActiveWindow
$benthicWindow.ScrollRow = 40
$benthicWindow.ScrollColumn = 1
$wsBenthic.Range("B61").Select() | Out-Null
